$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 1900
$ws.Range("B11").Value = 2515

$ws.Range("F21").Value = 170
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 2
$ws.Range("M21").Value = 2
$ws.Range("N21").Value = 2
$ws.Range("O21").Value = 2
$ws.Range("P21").Value = 2
$ws.Range("Q21").Value = 2
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 2
$ws.Range("T21").Value = 2

$ws.Range("F22").Value = 170
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 2
$ws.Range("N22").Value = 2
$ws.Range("O22").Value = 2
$ws.Range("P22").Value = 2
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 2
$ws.Range("T22").Value = 2

$ws.Range("F23").Value = 1105
$ws.Range("K23").Value = 2
$ws.Range("L23").Value = 2
$ws.Range("M23").Value = 2
$ws.Range("N23").Value = 2
$ws.Range("O23").Value = 2
$ws.Range("P23").Value = 2
$ws.Range("Q23").Value = 4
$ws.Range("R23").Value = 2
$ws.Range("S23").Value = 2
$ws.Range("T23").Value = 2

$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 6
$ws.Range("U24").Select() | Out-Null
